$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66; this shifts the existing rows 66-190
# down to 67-191 and extends the used range to A1:R191.
$ws.Rows(66).Insert()

# Populate the newly inserted row 66 with the new price record.
$ws.Range("A66").Value = 5
$ws.Range("B66").Value = "Macroferia Regional de Talca"
$ws.Range("C66").Value = "Maule"
$ws.Range("D66").Value = 44469
$ws.Range("E66").Value = 7
$ws.Range("F66").Value = 100114013
$ws.Range("G66").Value = "Zanahoria"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 400
$ws.Range("K66").Value = 7000
$ws.Range("L66").Value = 7000
$ws.Range("M66").Value = 7000
$ws.Range("N66").Value = "$/saco 20 kilos"
$ws.Range("O66").Value = "Región de Ñuble"
$ws.Range("P66").Value = 350
$ws.Range("Q66").Value = 20
$ws.Range("R66").Value = "Hortaliza"
